$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (24) with the latest gold price data.
$ws.Range("A24").Value = "15-10-2025"
$ws.Range("B24").Value = "The price of gold in India today is ₹12,889 per gram for 24 karat gold, ₹11,815 per gram for 22 karat gold and ₹9,697 per gram for 18 karat gold (also called 999 gold)."

# Match the formatting of the previous data row (border style + wrap text).
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("B23").Copy()
$ws.Range("B24").PasteSpecial(-4122)
